$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '91.918.00'
Set-TextValue $ws.Range('E2') '  -6.63%  '
Set-TextValue $ws.Range('D3') '3.300.51'
Set-TextValue $ws.Range('E3') '  -5.47%  '
Set-TextValue $ws.Range('E4') '  +0.14%  '
Set-TextValue $ws.Range('D5') '223.91'
Set-TextValue $ws.Range('E5') '  -11.58%  '
Set-TextValue $ws.Range('D6') '612.48'
Set-TextValue $ws.Range('E6') '  -8.10%  '
Set-TextValue $ws.Range('D7') '1.30'
Set-TextValue $ws.Range('E7') '  -12.34%  '
Set-TextValue $ws.Range('D8') '0.370'
Set-TextValue $ws.Range('E8') '  -13.56%  '
Set-TextValue $ws.Range('E9') '  +0.23%  '
Set-TextValue $ws.Range('D10') '0.890'
Set-TextValue $ws.Range('E10') '  -15.92%  '
Set-TextValue $ws.Range('D11') '3.292.96'
Set-TextValue $ws.Range('E11') '  -5.66%  '
Set-TextValue $ws.Range('D12') '0.189'
Set-TextValue $ws.Range('E12') '  -10.38%  '
Set-TextValue $ws.Range('D13') '38.81'
Set-TextValue $ws.Range('E13') '  -15.10%  '
Set-TextValue $ws.Range('D14') '91.892.92'
Set-TextValue $ws.Range('E14') '  -6.46%  '
Set-TextValue $ws.Range('D15') '5.73'
Set-TextValue $ws.Range('E15') '  -8.62%  '
Set-TextValue $ws.Range('D16') '3.919.37'
Set-TextValue $ws.Range('E16') '  -5.59%  '
Set-TextValue $ws.Range('D17') '0.0000236'
Set-TextValue $ws.Range('E17') '  -9.64%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D18') '3.297.39'
Set-TextValue $ws.Range('E18') '  -5.34%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D19') '7.63'
Set-TextValue $ws.Range('E19') '  -15.01%  '
Set-TextValue $ws.Range('D20') '16.32'
Set-TextValue $ws.Range('E20') '  -13.08%  '
Set-TextValue $ws.Range('D21') '10.57'
Set-TextValue $ws.Range('E21') '  -10.07%  '
Set-TextValue $ws.Range('D22') '477.70'
Set-TextValue $ws.Range('E22') '  -8.98%  '
Set-TextValue $ws.Range('D23') '0.426'
Set-TextValue $ws.Range('E23') '  -19.37%  '
Set-TextValue $ws.Range('D24') '3.02'
Set-TextValue $ws.Range('E24') '  -11.97%  '
Set-TextValue $ws.Range('D25') '0.0000177'
Set-TextValue $ws.Range('E25') '  -12.98%  '
Set-TextValue $ws.Range('D26') '5.96'
Set-TextValue $ws.Range('E26') '  -12.46%  '
Set-TextValue $ws.Range('D27') '87.88'
Set-TextValue $ws.Range('E27') '  -10.81%  '
Set-TextValue $ws.Range('D28') '3.502.67'
Set-TextValue $ws.Range('E28') '  -4.90%  '
Set-TextValue $ws.Range('D29') '11.09'
Set-TextValue $ws.Range('E29') '  -13.14%  '
Set-TextValue $ws.Range('E30') '  +0.54%  '
Set-TextValue $ws.Range('D31') '10.87'
Set-TextValue $ws.Range('E31') '  -12.74%  '
Set-TextValue $ws.Range('D32') '2.54'
Set-TextValue $ws.Range('E32') '  -11.52%  '
$ws.Range('B33').Value = 'Binance-PegBSC-USD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range('D33') '1.00'
Set-TextValue $ws.Range('E33') '  -0.32%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D34') '0.127'
Set-TextValue $ws.Range('E34') '  -12.48%  '
Set-TextValue $ws.Range('D35') '0.166'
Set-TextValue $ws.Range('E35') '  -12.83%  '
Set-TextValue $ws.Range('D36') '27.68'
Set-TextValue $ws.Range('E36') '  -11.51%  '
Set-TextValue $ws.Range('D37') '0.507'
Set-TextValue $ws.Range('E37') '  -15.27%  '
Set-TextValue $ws.Range('E38') '  -0.03%  '
Set-TextValue $ws.Range('D39') '508.67'
Set-TextValue $ws.Range('E39') '  -3.53%  '
Set-TextValue $ws.Range('D40') '7.17'
Set-TextValue $ws.Range('E40') '  -10.01%  '
Set-TextValue $ws.Range('D41') '1.34'
Set-TextValue $ws.Range('E41') '  -11.55%  '
Set-TextValue $ws.Range('D42') '0.143'
Set-TextValue $ws.Range('E42') '  -8.55%  '
Set-TextValue $ws.Range('D43') '0.847'
Set-TextValue $ws.Range('E43') '  -7.87%  '
Set-TextValue $ws.Range('D44') '23.99'
Set-TextValue $ws.Range('E44') '  -1.77%  '
Set-TextValue $ws.Range('D45') '1.62'
Set-TextValue $ws.Range('E45') '  -8.55%  '
Set-TextValue $ws.Range('D46') '3.48'
Set-TextValue $ws.Range('E46') '  -4.38%  '
Set-TextValue $ws.Range('D47') '5.24'
Set-TextValue $ws.Range('E47') '  -9.62%  '
Set-TextValue $ws.Range('D48') '2.06'
Set-TextValue $ws.Range('E48') '  -7.66%  '
Set-TextValue $ws.Range('D49') '51.84'
Set-TextValue $ws.Range('E49') '  -6.74%  '
Set-TextValue $ws.Range('D50') '0.0379'
Set-TextValue $ws.Range('E50') '  -12.91%  '
Set-TextValue $ws.Range('D51') '2.99'
Set-TextValue $ws.Range('E51') '  -7.93%  '
